# Updates cryptos list figures (price & volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.684.69'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.585.85'
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("E4").Value = '  +1.33%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '206.27'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("E6").Value = '  -1.93%  '
$ws.Range("E7").Value = '  +1.37%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '22.24'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("E11").Value = '  -0.92%  '
$ws.Range("D12").Value = '1.810.30'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("D13").Value = '1.590.91'
$ws.Range("E13").Value = '  -1.66%  '
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("E15").Value = '  -4.44%  '
$ws.Range("D16").Value = '27.652.40'
$ws.Range("E16").Value = '  -0.05%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '63.30'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -2.16%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '219.53'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -3.46%  '
$ws.Range("E19").Value = '  -2.94%  '
$ws.Range("E20").Value = '  -4.02%  '
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("E23").Value = '  -4.62%  '
$ws.Range("E24").Value = '  -2.18%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '155.32'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.49%  '
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("E27").Value = '  +1.33%  '
$ws.Range("E28").Value = '  -2.04%  '
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("E30").Value = '  -1.41%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.0466'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("D33").Value = '1.380.62'
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("E34").Value = '  -4.30%  '
$ws.Range("E35").Value = '  -3.15%  '
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("E38").Value = '  -2.61%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.538'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -2.95%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.823'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -2.14%  '
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("E42").Value = '  -2.49%  '
$ws.Range("E43").Value = '  -2.84%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '63.52'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.98%  '
$ws.Range("E45").Value = '  +2.76%  '
$ws.Range("E46").Value = '  -2.57%  '
$ws.Range("D47").Value = '1.721.32'
$ws.Range("E47").Value = '  -1.99%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '88.37'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("E49").Value = '  +13.21%  '
$ws.Range("E50").Value = '  -3.55%  '
$ws.Range("E51").Value = '  -0.51%  '
